$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (label "0.0")
$ws.Range("B2").Value = 0.9748502994011976
$ws.Range("C2").Value = 0.8357289527720739
$ws.Range("D2").Value = 0.8999447208402431
$ws.Range("E2").Value = 974

# Row 3 (label "1.0")
$ws.Range("B3").Value = 0.3360995850622407
$ws.Range("C3").Value = 0.7941176470588235
$ws.Range("D3").Value = 0.4723032069970846
$ws.Range("E3").Value = 102

# Row 4 (accuracy)
$ws.Range("B4").Value = 0.8317843866171004
$ws.Range("C4").Value = 0.8317843866171004
$ws.Range("D4").Value = 0.8317843866171004
$ws.Range("E4").Value = 0.8317843866171004

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.6554749422317192
$ws.Range("C5").Value = 0.8149232999154488
$ws.Range("D5").Value = 0.6861239639186638
$ws.Range("E5").Value = 1076

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.9142995811274303
$ws.Range("C6").Value = 0.8317843866171004
$ws.Range("D6").Value = 0.8594062130223973
$ws.Range("E6").Value = 1076
